$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C31").Value = 50
$ws.Range("E31").Value = 0.02162629757785467

$ws.Range("C36").Value = 135
$ws.Range("E36").Value = 0.06994818652849741

$ws.Range("C37").Value = 861
$ws.Range("D37").Value = 861
